# The AlcoholUse "StartDate" mapping (columns B..J: zib, name_zib, path_zib,
# alias_zib, type_zib, card._zib, stereotype_zib, id_zib, definition_zib) was
# attached to the wrong xtehr row (EHDSSubstanceUse.period, row 13). It
# belongs on EHDSSubstanceUse.frequencyAndQuantity.period (row 16) instead.
# Move (cut/paste) the block of values from row 13 to row 16.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$src = $ws.Range("B13:J13")
$dst = $ws.Range("B16:J16")

# Grab the values with Value2 (not the buggy chained "Value" getter) so we can
# relocate them safely even though src and dst overlap other writes.
$moved = $src.Value2

$dst.Value2 = $moved
$src.Value2 = ""
